# Scheduled market-data refresh: update cached Universalis price/profit
# columns (H:N) per leve row across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1467.375
$ws.Range("I4").Value = 1431
$ws.Range("K4").Value = 1431
$ws.Range("M4").Value = -1317
# Row 69
$ws.Range("H69").Value = 11249.25
$ws.Range("J69").Value = 11249.25
$ws.Range("L69").Value = 33747.75
$ws.Range("N69").Value = -35495.75
# Row 70
$ws.Range("H70").Value = 126017240
$ws.Range("I70").Value = 253049300
$ws.Range("J70").Value = 62501212
$ws.Range("K70").Value = 759147900
$ws.Range("L70").Value = 187503636
$ws.Range("M70").Value = -759147630
$ws.Range("N70").Value = -187504176
# Row 72
$ws.Range("H72").Value = 11249.25
$ws.Range("J72").Value = 11249.25
$ws.Range("L72").Value = 101243.25
$ws.Range("N72").Value = -109979.25
# Row 73
$ws.Range("H73").Value = 126017240
$ws.Range("I73").Value = 253049300
$ws.Range("J73").Value = 62501212
$ws.Range("K73").Value = 759147900
$ws.Range("L73").Value = 187503636
$ws.Range("M73").Value = -759146964
$ws.Range("N73").Value = -187505508
# Row 86
$ws.Range("H86").Value = 10026
$ws.Range("I86").Value = 4420.5
$ws.Range("K86").Value = 4420.5
$ws.Range("M86").Value = -3297.5
# Row 89
$ws.Range("H89").Value = 10026
$ws.Range("I89").Value = 4420.5
$ws.Range("K89").Value = 22102.5
$ws.Range("M89").Value = -16486.5
$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 132
$ws.Range("H132").Value = 30005004
$ws.Range("I132").Value = 5846.7144
$ws.Range("J132").Value = 100003030
$ws.Range("K132").Value = 17540.1432
$ws.Range("L132").Value = 300009090
$ws.Range("M132").Value = -15010.1432
$ws.Range("N132").Value = -300014150
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1801.5
$ws.Range("I94").Value = 1850.3334
$ws.Range("K94").Value = 1850.3334
$ws.Range("M94").Value = -1399.3334
# Row 99
$ws.Range("H99").Value = 2549.15
$ws.Range("I99").Value = 2610.1667
$ws.Range("K99").Value = 2610.1667
$ws.Range("M99").Value = -1112.1667
$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1580.2222
$ws.Range("J68").Value = 1993.25
$ws.Range("L68").Value = 5979.75
$ws.Range("N68").Value = -7601.75
# Row 71
$ws.Range("H71").Value = 1580.2222
$ws.Range("J71").Value = 1993.25
$ws.Range("L71").Value = 17939.25
$ws.Range("N71").Value = -26051.25
# Row 140
$ws.Range("H140").Value = 7334.8667
$ws.Range("I140").Value = 2865.8
$ws.Range("K140").Value = 8597.400000000001
$ws.Range("M140").Value = -3417.400000000001
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 25
$ws.Range("J2").Value = 25
$ws.Range("L2").Value = 25
$ws.Range("N2").Value = -251
# Row 97
$ws.Range("H97").Value = 931.5769
$ws.Range("I97").Value = 607.4737
$ws.Range("J97").Value = 1811.2858
$ws.Range("K97").Value = 607.4737
$ws.Range("L97").Value = 1811.2858
$ws.Range("M97").Value = -111.4737
$ws.Range("N97").Value = -2803.2858
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3084.8
$ws.Range("I40").Value = 3182.2856
$ws.Range("J40").Value = 2857.3333
$ws.Range("K40").Value = 3182.2856
$ws.Range("L40").Value = 2857.3333
$ws.Range("M40").Value = -3046.2856
$ws.Range("N40").Value = -3129.3333
# Row 61
$ws.Range("H61").Value = 76925710
$ws.Range("I61").Value = 83335950
$ws.Range("J61").Value = 2795
$ws.Range("K61").Value = 83335950
$ws.Range("L61").Value = 2795
$ws.Range("M61").Value = -83335748
$ws.Range("N61").Value = -3199
# Row 93
$ws.Range("H93").Value = 3477503.5
$ws.Range("I93").Value = 1812.5
$ws.Range("K93").Value = 1812.5
$ws.Range("M93").Value = -564.5
# Row 100
$ws.Range("H100").Value = 13174720
$ws.Range("I100").Value = 3903.125
$ws.Range("J100").Value = 22753496
$ws.Range("K100").Value = 3903.125
$ws.Range("L100").Value = 22753496
$ws.Range("M100").Value = -3362.125
$ws.Range("N100").Value = -22754578
# Row 113
$ws.Range("H113").Value = 76925710
$ws.Range("I113").Value = 83335950
$ws.Range("J113").Value = 2795
$ws.Range("K113").Value = 83335950
$ws.Range("L113").Value = 2795
$ws.Range("M113").Value = -83333780
$ws.Range("N113").Value = -7135
# Row 132
$ws.Range("H132").Value = 4670.65
$ws.Range("I132").Value = 3271
$ws.Range("J132").Value = 5815.8184
$ws.Range("K132").Value = 9813
$ws.Range("L132").Value = 17447.4552
$ws.Range("M132").Value = -7283
$ws.Range("N132").Value = -22507.4552
# Row 136
$ws.Range("H136").Value = 2331.1
$ws.Range("I136").Value = 2171.077
$ws.Range("K136").Value = 6513.231000000001
$ws.Range("M136").Value = -3963.231000000001
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2031934.1
$ws.Range("I62").Value = 8600
$ws.Range("J62").Value = 3043601.2
$ws.Range("K62").Value = 8600
$ws.Range("L62").Value = 3043601.2
$ws.Range("M62").Value = -7976
$ws.Range("N62").Value = -3044849.2
# Row 65
$ws.Range("H65").Value = 2031934.1
$ws.Range("I65").Value = 8600
$ws.Range("J65").Value = 3043601.2
$ws.Range("K65").Value = 43000
$ws.Range("L65").Value = 15218006
$ws.Range("M65").Value = -39880
$ws.Range("N65").Value = -15224246
# Row 113
$ws.Range("H113").Value = 511.8
$ws.Range("I113").Value = 604.2857
$ws.Range("J113").Value = 296
$ws.Range("K113").Value = 1812.8571
$ws.Range("L113").Value = 888
$ws.Range("M113").Value = 357.1428999999998
$ws.Range("N113").Value = -5228
# Row 122
$ws.Range("H122").Value = 2768.6667
$ws.Range("I122").Value = 2819.8096
$ws.Range("K122").Value = 8459.4288
$ws.Range("M122").Value = -6009.4288
# Row 126
$ws.Range("H126").Value = 33336786
$ws.Range("I126").Value = 43481612
$ws.Range("J126").Value = 3784.1428
$ws.Range("K126").Value = 130444836
$ws.Range("L126").Value = 11352.4284
$ws.Range("M126").Value = -130442366
$ws.Range("N126").Value = -16292.4284
# Row 132
$ws.Range("H132").Value = 347226.34
$ws.Range("I132").Value = 2083.762
$ws.Range("K132").Value = 6251.286
$ws.Range("M132").Value = -3721.286
# Row 136
$ws.Range("H136").Value = 160277.84
$ws.Range("I136").Value = 6855.68
$ws.Range("J136").Value = 671685.0600000001
$ws.Range("K136").Value = 20567.04
$ws.Range("L136").Value = 2015055.18
$ws.Range("M136").Value = -18017.04
$ws.Range("N136").Value = -2020155.18
